# Applies the "new run with new features" update to the ordinal regression
# evaluation worksheet.
#
# Summary of change:
#  - Row 2 label changes from "LogisticAT" to "LAD"
#  - Row 3 label changes from "LAD" to "LogisticAT"
#    (i.e. the two model names swap places)
#  - Rows 2-5, columns B-E get new metric values from the re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels (column A) ---
$ws.Range("A2").Value = "LAD"
$ws.Range("A3").Value = "LogisticAT"

# --- Row 2 metrics ---
$ws.Range("B2").Value = 0.4808
$ws.Range("C2").Value = 0.4808
$ws.Range("D2").Value = 0.9121
$ws.Range("E2").Value = 0.614

# --- Row 3 metrics ---
$ws.Range("B3").Value = 0.4615
$ws.Range("C3").Value = 0.4615
$ws.Range("D3").Value = 0.9121
$ws.Range("E3").Value = 0.6304999999999999

# --- Row 4 metrics ---
$ws.Range("B4").Value = 0.456
$ws.Range("C4").Value = 0.456
$ws.Range("D4").Value = 0.9121
$ws.Range("E4").Value = 0.636

# --- Row 5 metrics ---
$ws.Range("B5").Value = 0.4945
$ws.Range("C5").Value = 0.4945
$ws.Range("D5").Value = 0.8915
$ws.Range("E5").Value = 0.6277
